$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '60.389.11'
$ws.Range("E2").Value = '  +3.16%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.329.89'
$ws.Range("E3").Value = '  +1.22%  '
$ws.Range("E4").Value = '  -0.02%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '544.65'
$ws.Range("E5").Value = '  +1.43%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '131.27'
$ws.Range("E6").Value = '  -0.60%  '
$ws.Range("E7").Value = '  -0.03%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.577'
$ws.Range("E8").Value = '  -1.40%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '2.326.95'
$ws.Range("E9").Value = '  +1.15%  '
$ws.Range("E10").Value = '  +0.82%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '5.53'
$ws.Range("E11").Value = '  +0.60%  '
$ws.Range("E12").Value = '  +0.18%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.334'
$ws.Range("E13").Value = '  +0.39%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '23.65'
$ws.Range("E14").Value = '  -0.24%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '2.745.75'
$ws.Range("E15").Value = '  +1.21%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '60.346.89'
$ws.Range("E16").Value = '  +3.21%  '
$ws.Range("E17").Value = '  +0.03%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '2.323.90'
$ws.Range("E18").Value = '  +1.03%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '10.59'
$ws.Range("E19").Value = '  +0.21%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '4.14'
$ws.Range("E20").Value = '  -1.34%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '315.06'
$ws.Range("E21").Value = '  -0.31%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.65'
$ws.Range("E22").Value = '  +0.94%  '
$ws.Range("E23").Value = '  -0.09%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '63.99'
$ws.Range("E24").Value = '  +1.40%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.172'
$ws.Range("E25").Value = '  +1.76%  '
$ws.Range("E26").Value = '  +0.09%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '7.85'
$ws.Range("E27").Value = '  -1.10%  '
$ws.Range("E28").Value = '  +5.01%  '
$ws.Range("E29").Value = '  +9.89%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '173.12'
$ws.Range("E30").Value = '  +0.91%  '
$ws.Range("E31").Value = '  +1.08%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.0₃0733'
$ws.Range("E32").Value = '  +1.05%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '5.94'
$ws.Range("E33").Value = '  +1.98%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.37'
$ws.Range("E34").Value = '  +10.13%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.380'
$ws.Range("E35").Value = '  -1.05%  '
$ws.Range("E36").Value = '  -0.01%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '17.84'
$ws.Range("E37").Value = '  -0.27%  '
$ws.Range("E38").Value = '  -0.02%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '4.07'
$ws.Range("E39").Value = '  +1.94%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '321.02'
$ws.Range("E40").Value = '  +10.83%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.53'
$ws.Range("E41").Value = '  +1.75%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '37.96'
$ws.Range("E42").Value = '  -0.99%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '137.39'
$ws.Range("E43").Value = '  -2.71%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '3.49'
$ws.Range("E44").Value = '  +1.08%  '
$ws.Range("E45").Value = '  -1.35%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '19.18'
$ws.Range("E46").Value = '  +4.95%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.563'
$ws.Range("E47").Value = '  +1.27%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.0495'
$ws.Range("E48").Value = '  +0.02%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.0213'
$ws.Range("E49").Value = '  +1.56%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0₆0215'
$ws.Range("E50").Value = '  +16.90%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '11.02'
$ws.Range("E51").Value = '  +0.61%  '
